$d = $word.ActiveDocument

# 1. Ativação date bump: 2020 -> 2023
$d.Content.Find.Execute(
    "Ativação: 01/01/2020", $true, $false, $false, $false, $false,
    $true, 1, $false, "Ativação: 01/01/2023", 2) | Out-Null

# 2. Objetivos paragraph rewrite
$d.Content.Find.Execute(
    "-Capacitar o aluno a escrever e balancear reações químicas, mostrando os produtos esperados, para os elementos da tabela periódica e seus compostos.-Capacitar o aluno para relacionar as propriedades químicas e físicas dos elementos e seus compostos com suas posições na tabela periódica.-Capacitar o aluno a escrever os métodos de obtenção dos elementos e seus compostos, bem como descrever suas aplicações.",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "- Capacitar o aluno para relacionar as propriedades químicas e físicas dos elementos e seus compostos com suas posições na tabela periódica.-Capacitar o aluno a escrever os métodos industriais de obtenção dos elementos e seus compostos, bem como descrever suas aplicações- capacitar o aluno a comunicar-se eficazmente nas formas escrita, oral e gráfica",
    2) | Out-Null

# 3. Insert a new "Ângelo Capri Neto" run+break before the existing
#    "5840963 - Daniela Camargo Vernilli" run in the Docente(s) list.
$r = $d.Content
$r.Find.Execute("5840963 - Daniela Camargo Vernilli", $true, $false, $false,
                 $false, $false, $true, 1, $false, "", 0) | Out-Null
$ins = $r.Duplicate
$ins.Collapse(1)
$ins.InsertBefore("5840712 - Ângelo Capri Neto" + [char]11)

# 4. Programa paragraph rewrite
$d.Content.Find.Execute(
    "- Metais Representativos: Características gerais dos metais dos Grupos 1, 2 e 13.- Metais de transição: Introdução e Propriedades gerais, Complexos.",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Metais e compostos dos grupos 1, 2, 13 e de transição da Tabela Periódica: Propriedades físicas e químicas (relação com a posição na Tabela Periódica), processos de obtenção dos metais e compostos e aplicações - Formação de Complexos.Relacionar a disciplina com disciplinas anteriores e posteriores da grade do curso.",
    2) | Out-Null

# 5. Método text rewrite
$d.Content.Find.Execute(
    "Duas provas bimestrais escritas (P1 e P2), cada uma valendo nota de 0,0  a 10,0.",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "A avaliação tem como requisito quantificar as competências adquiridas conforme objetivadas.Duas provas escritas (P1 e P2) e listas de exercícios de acompanhamento continuado. A partir das notas das listas de exercício será calculada a média, LE.",
    2) | Out-Null

# 6. Critério text rewrite
$d.Content.Find.Execute(
    "MS= P1+P2/2, onde: MS= média do semestre.MS> ou = 5,0 = Aluno AprovadoMS< 3,0 = Aluno Reprovado3,0 < ou = MS < 5,0 = Aluno de Recuperação.",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "NF = (P1 + P2 + LE) /3",
    2) | Out-Null

# 7. Norma de recuperação text rewrite
$d.Content.Find.Execute(
    "Estudo dirigido de todo o conteúdo da disciplina e uma prova (PR) valendo nota de 0,0 a 10,0, contendo todo o conteúdo da disciplina.O aluno será aprovado se apresentar (média final) MF > ou = 5,0.Onde: MF= MS+PR/2, onde:  MS= média do semestre e PR= prova de recuperação.",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Será realizada uma prova escrita valendo de zero a dez (NR) e a média final calculada pela equação: NF + NR",
    2) | Out-Null

# 8. Bibliografia text rewrite
$d.Content.Find.Execute(
    "- LEE, J. D. “Química Inorgânica não tão Concisa”, Editora Edgard Blücher, 1999.- SHRIVER, D. F.; ATKINS, P. W. “Química Inorgânica”, Editora Bookman, 4ª edição, 2008.- QUAGLIANO, J. V.; VALLARINO, L. “Química”, Editora Guanabara Koogan, 1973.- MELLOR, J. W. “Química Inorgânica Moderna”, Editora: Globo – Porto Alegre, 1967.- GREENWOOD, N. N.; EARNSHAW, A. “Chemistry of the Elements”, Butterworth Heinemann, 1997.- BUCHEL, K. H.; MORETTO, H. H.; WODITSCH, P. “Industrial Inorganic Chemistry”, Editora Wiley-VCH, 2000.- RAYNER-CANHAM, G.; OVERTON, T. “Química Inorgânica Descritiva”, Editora: Gen-LTC, 5ª edição, 2015.- SOUZA, M.M.V.M. “Processos Inorgânicos”, Editora: Synergia, 1ª edição, 2012.",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "LEE, J. D. “Química Inorgânica não tão Concisa”, Editora Edgard Blücher, 1999. - SHRIVER, D. F.; ATKINS, P. W. “Química Inorgânica”, Editora Bookman, 4ª edição, 2008. - QUAGLIANO, J. V.; VALLARINO, L. “Química”, Editora Guanabara Koogan, 1973. - BUCHEL, K. H.; MORETTO, H. H.; WODITSCH, P. “Industrial Inorganic Chemistry”, Editora Wiley-VCH, 2000. - RAYNER-CANHAM, G.; OVERTON, T. “Química Inorgânica Descritiva”, Editora: Gen-LTC, 5ª edição, 2015. - SOUZA, M.M.V.M. “Processos Inorgânicos”, Editora: Synergia, 1ª edição, 2012.",
    2) | Out-Null

Write-Output "done"
